# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that is safe as literal text (never parses as a pure number).
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Helper: write a numeric-looking string while keeping the cell a plain text/string
# cell (no inline-formula, no NumberFormat change that would bake a new style).
# Route it through a temporary text-formula + paste-special-values so Excel
# stores the literal characters instead of coercing them to a number.
function Set-LiteralText($addr, $val) {
    $escaped = $val.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

Set-TextValue 'D2' '67.183.10'

Set-TextValue 'D3' '2.473.87'
Set-TextValue 'E3' '  -2.17%  '

Set-LiteralText 'D4' '1.00'
Set-TextValue 'E4' '  +0.02%  '

Set-LiteralText 'D5' '583.17'
Set-TextValue 'E5' '  -1.39%  '

Set-LiteralText 'D6' '169.05'
Set-TextValue 'E6' '  -1.35%  '

Set-TextValue 'E7' '  +0.05%  '

Set-TextValue 'E8' '  -1.99%  '

Set-TextValue 'D9' '2.475.70'

Set-TextValue 'E10' '  -2.44%  '

Set-TextValue 'E12' '  -2.63%  '

Set-TextValue 'E13' '  -3.11%  '

Set-TextValue 'E14' '  -3.07%  '

Set-TextValue 'E15' '  -1.07%  '

Set-TextValue 'D16' '66.771.35'
Set-TextValue 'E16' '  -1.04%  '

Set-TextValue 'E17' '  -4.02%  '

Set-TextValue 'D18' '2.446.15'
Set-TextValue 'E18' '  -2.04%  '

Set-TextValue 'E19' '  -5.34%  '

Set-TextValue 'E20' '  -3.37%  '

Set-LiteralText 'D21' '354.49'
Set-TextValue 'E21' '  -3.76%  '

Set-TextValue 'E23' '  -0.60%  '

Set-LiteralText 'D24' '69.18'
Set-TextValue 'E24' '  -3.47%  '

Set-LiteralText 'D25' '4.24'
Set-TextValue 'E25' '  -7.51%  '

Set-TextValue 'E26' '  -6.60%  '

Set-TextValue 'E27' '  -6.87%  '

Set-LiteralText 'D28' '0.999'
Set-TextValue 'E28' '  +0.05%  '

Set-TextValue 'D29' '2.591.88'
Set-TextValue 'E29' '  -1.86%  '

Set-TextValue 'D30' '0.0₃0909'
Set-TextValue 'E30' '  -5.59%  '

Set-LiteralText 'D31' '518.04'
Set-TextValue 'E31' '  -3.67%  '

Set-LiteralText 'D32' '7.76'
Set-TextValue 'E32' '  -7.07%  '

Set-TextValue 'E33' '  -5.46%  '

Set-TextValue 'E34' '  -5.15%  '

Set-LiteralText 'D35' '1.00'
Set-TextValue 'E35' '  +0.08%  '

Set-LiteralText 'D36' '0.120'
Set-TextValue 'E36' '  -6.91%  '

Set-LiteralText 'D37' '157.90'
Set-TextValue 'E37' '  -0.43%  '

Set-LiteralText 'D38' '18.68'
Set-TextValue 'E38' '  +0.25%  '

Set-LiteralText 'D39' '18.42'
Set-TextValue 'E39' '  -3.76%  '

Set-TextValue 'E40' '  -5.09%  '

Set-TextValue 'E41' '  -0.02%  '

Set-LiteralText 'D42' '0.327'
Set-TextValue 'E42' '  -6.63%  '

Set-TextValue 'B43' 'Stacks'
Set-TextValue 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-LiteralText 'D43' '1.67'
Set-TextValue 'E43' '  -6.21%  '

Set-TextValue 'B44' 'RenderToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-LiteralText 'D44' '4.80'
Set-TextValue 'E44' '  -6.39%  '

Set-TextValue 'E45' '  -4.60%  '

Set-TextValue 'E46' '  -2.35%  '

Set-LiteralText 'D47' '141.35'
Set-TextValue 'E47' '  -3.50%  '

Set-TextValue 'E48' '  -6.39%  '

Set-TextValue 'E49' '  -6.59%  '

Set-TextValue 'E50' '  -10.60%  '

Set-TextValue 'E51' '  -7.29%  '

$excel.CutCopyMode = 0
